# Apply workbook updates for "4.0 files and mdl" refresh
$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the "last updated" date stamp (C1) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45392

# --- "MCF" sheet: update maximum capacity factor assumptions to 1 ---
$wsMcf = $wb.Worksheets.Item("MCF")

$wsMcf.Range("B2").Value = 1    # hard coal
$wsMcf.Range("B3").Value = 1    # natural gas steam turbine
$wsMcf.Range("B4").Value = 1    # natural gas combined cycle
$wsMcf.Range("B6").Value = 1    # hydro
$wsMcf.Range("B10").Value = 1   # biomass
$wsMcf.Range("B11").Value = 1   # geothermal
$wsMcf.Range("B12").Value = 1   # petroleum
$wsMcf.Range("B13").Value = 1   # natural gas peaker
$wsMcf.Range("B14").Value = 1   # lignite
$wsMcf.Range("B16").Value = 1   # crude oil
$wsMcf.Range("B17").Value = 1   # heavy or residual fuel oil
$wsMcf.Range("B18").Value = 1   # municipal solid waste

# B19:B25 are formulas referencing the cells above; they recalc automatically.

# --- Restore the active selection on the MCF sheet (now on B17) ---
$wsMcf.Activate()
$wsMcf.Range("B17").Select()
